$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Participant" tab query (B2) with the new, more elaborate
# Cypher query that matches participants independently of their samples
# (so participants without samples/files still show up), then re-derives
# the related samples via a second OPTIONAL MATCH pass.
$newParticipantQuery = "MATCH (p:participant)-->(s:study)`n" + `
  "OPTIONAL MATCH (samp:sample)-->(p)`n" + `
  "OPTIONAL MATCH (p)<--(diag:diagnosis)`n" + `
  "OPTIONAL MATCH (samp)<--(f:file)`n" + `
  "OPTIONAL MATCH (f)<--(g:genomic_info)`n" + `
  "WITH s, p, samp, f, g, diag`n" + `
  "WHERE g.platform in ['Illumina Next Seq 500']`n" + `
  "with p`n" + `
  "OPTIONAL MATCH (p)-->(s:study)`n" + `
  "OPTIONAL MATCH (samp:sample)-->(p)`n" + `
  "WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp`n" + `
  "RETURN`n" + `
  "coalesce(p.participant_id,'') as ``Participant ID``,`n" + `
  "coalesce(s.study_name, '') as ``Study Name``,`n" + `
  "coalesce(s.phs_accession,'') as ``Accession``,`n" + `
  "coalesce(p.gender,'') as ``Gender``,`n" + `
  "coalesce(apoc.text.join(samp, ','), '') as ``Samples```n" + `
  "ORDER BY p.participant_id LIMIT 100"

$ws.Range("B2").Value = $newParticipantQuery

# The new query text wraps over more lines, so the row grows taller.
$ws.Rows("2:2").RowHeight = 279

# The active selection/scroll position moved when the sheet was last saved.
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("B5").Select()
